# Updates the cryptocurrency price/volume/date-flag data in the worksheet
# to match the latest scrape, per commit:
# "Updated symbol list on Mon Jan 30 01:23:59 UTC 2023 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'317.31"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'2.61%"
$ws.Range("E2").Style = "Normal"
$ws.Range("G2").Value = "'1"
$ws.Range("G2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'39.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'0.39%"
$ws.Range("E3").Style = "Normal"
$ws.Range("G3").Value = "'1"
$ws.Range("G3").Style = "Normal"

# Row 4
$ws.Range("D4").Value = "'5.125"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.57%"
$ws.Range("E4").Style = "Normal"
$ws.Range("G4").Value = "'1"
$ws.Range("G4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'0.08198"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'0.36%"
$ws.Range("E5").Style = "Normal"
$ws.Range("G5").Value = "'1"
$ws.Range("G5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'2.017"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'3.38%"
$ws.Range("E6").Style = "Normal"
$ws.Range("G6").Value = "'1"
$ws.Range("G6").Style = "Normal"

# Row 7
$ws.Range("B7").Value = "'GateToken"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = "'4.335"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'2.70%"
$ws.Range("E7").Style = "Normal"
$ws.Range("G7").Value = "'1"
$ws.Range("G7").Style = "Normal"

# Row 8
$ws.Range("B8").Value = "'KuCoinToken"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = "'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = "'8.322"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'2.73%"
$ws.Range("E8").Style = "Normal"
$ws.Range("G8").Value = "'1"
$ws.Range("G8").Style = "Normal"

# Row 9
$ws.Range("B9").Value = "'MXToken"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = "'0.9418"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.95%"
$ws.Range("E9").Style = "Normal"
$ws.Range("G9").Value = "'1"
$ws.Range("G9").Style = "Normal"

# Row 10
$ws.Range("B10").Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'0.1288"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-6.07%"
$ws.Range("E10").Style = "Normal"
$ws.Range("G10").Value = "'1"
$ws.Range("G10").Style = "Normal"

# Row 11
$ws.Range("B11").Value = "'WazirX"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'0.1992"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'4.10%"
$ws.Range("E11").Style = "Normal"
$ws.Range("G11").Value = "'1"
$ws.Range("G11").Style = "Normal"

# Row 12
$ws.Range("B12").Value = "'MandalaExchangeToken"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'0.09045"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-0.19%"
$ws.Range("E12").Style = "Normal"
$ws.Range("G12").Value = "'1"
$ws.Range("G12").Style = "Normal"

# Row 13
$ws.Range("B13").Value = "'BitrueCoin"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'0.03512"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.47%"
$ws.Range("E13").Style = "Normal"
$ws.Range("G13").Value = "'1"
$ws.Range("G13").Style = "Normal"

# Row 14
$ws.Range("B14").Value = "'BitMartToken"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'0.09784"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.50%"
$ws.Range("E14").Style = "Normal"
$ws.Range("G14").Value = "'1"
$ws.Range("G14").Style = "Normal"

# Row 15
$ws.Range("B15").Value = "'BitForexToken"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = "'0.001412"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'1.70%"
$ws.Range("E15").Style = "Normal"
$ws.Range("G15").Value = "'1"
$ws.Range("G15").Style = "Normal"

# Row 16
$ws.Range("B16").Value = "'TigerCash"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'0.005999"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'2.42%"
$ws.Range("E16").Style = "Normal"
$ws.Range("G16").Value = "'1"
$ws.Range("G16").Style = "Normal"

# Row 17
$ws.Range("B17").Value = "'LEO"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'3.693"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-2.79%"
$ws.Range("E17").Style = "Normal"
$ws.Range("G17").Value = "'1"
$ws.Range("G17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'3.304"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-1.45%"
$ws.Range("E18").Style = "Normal"
$ws.Range("G18").Value = "'1"
$ws.Range("G18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'0.3489"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.41%"
$ws.Range("E19").Style = "Normal"
$ws.Range("G19").Value = "'1"
$ws.Range("G19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'0.1322"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-1.26%"
$ws.Range("E20").Style = "Normal"
$ws.Range("G20").Value = "'1"
$ws.Range("G20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'4.974"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'7.03%"
$ws.Range("E21").Style = "Normal"
$ws.Range("G21").Value = "'1"
$ws.Range("G21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'0.2506"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'4.10%"
$ws.Range("E22").Style = "Normal"
$ws.Range("G22").Value = "'1"
$ws.Range("G22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'0.04352"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.43%"
$ws.Range("E23").Style = "Normal"
$ws.Range("G23").Value = "'1"
$ws.Range("G23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'0.001242"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'3.12%"
$ws.Range("E24").Style = "Normal"
$ws.Range("G24").Value = "'1"
$ws.Range("G24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'0.004767"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'8.83%"
$ws.Range("E25").Style = "Normal"
$ws.Range("G25").Value = "'1"
$ws.Range("G25").Style = "Normal"

# Row 26
$ws.Range("E26").Value = "'0.77%"
$ws.Range("E26").Style = "Normal"
$ws.Range("G26").Value = "'1"
$ws.Range("G26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "'0.0003699"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'-6.98%"
$ws.Range("E27").Style = "Normal"
$ws.Range("G27").Value = "'1"
$ws.Range("G27").Style = "Normal"

# Row 28
$ws.Range("G28").Value = "'1"
$ws.Range("G28").Style = "Normal"

# Row 29
$ws.Range("G29").Value = "'1"
$ws.Range("G29").Style = "Normal"

# Row 30
$ws.Range("G30").Value = "'1"
$ws.Range("G30").Style = "Normal"

# Row 31
$ws.Range("G31").Value = "'1"
$ws.Range("G31").Style = "Normal"

# Row 32
$ws.Range("G32").Value = "'1"
$ws.Range("G32").Style = "Normal"

# Row 33
$ws.Range("G33").Value = "'1"
$ws.Range("G33").Style = "Normal"

# Row 34
$ws.Range("G34").Value = "'1"
$ws.Range("G34").Style = "Normal"

# Row 35
$ws.Range("G35").Value = "'1"
$ws.Range("G35").Style = "Normal"

# Row 36
$ws.Range("G36").Value = "'1"
$ws.Range("G36").Style = "Normal"

# Row 37
$ws.Range("G37").Value = "'1"
$ws.Range("G37").Style = "Normal"

# Row 38
$ws.Range("G38").Value = "'1"
$ws.Range("G38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'0.02265"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'10.00%"
$ws.Range("E39").Style = "Normal"
$ws.Range("G39").Value = "'1"
$ws.Range("G39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'0.05187"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'2.59%"
$ws.Range("E40").Style = "Normal"
$ws.Range("G40").Value = "'1"
$ws.Range("G40").Style = "Normal"

# Row 41
$ws.Range("D41").Value = "'0.007757"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'5.65%"
$ws.Range("E41").Style = "Normal"
$ws.Range("G41").Value = "'1"
$ws.Range("G41").Style = "Normal"

# Row 42
$ws.Range("D42").Value = "'0.01049"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'7.81%"
$ws.Range("E42").Style = "Normal"
$ws.Range("G42").Value = "'1"
$ws.Range("G42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = "'0.1405"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'2.39%"
$ws.Range("E43").Style = "Normal"
$ws.Range("G43").Value = "'1"
$ws.Range("G43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'0.002048"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-3.25%"
$ws.Range("E44").Style = "Normal"
$ws.Range("G44").Value = "'1"
$ws.Range("G44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'0.009040"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-0.91%"
$ws.Range("E45").Style = "Normal"
$ws.Range("G45").Value = "'1"
$ws.Range("G45").Style = "Normal"

# Row 46
$ws.Range("D46").Value = "'0.00006910"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'12.27%"
$ws.Range("E46").Style = "Normal"
$ws.Range("G46").Value = "'1"
$ws.Range("G46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.88%"
$ws.Range("E47").Style = "Normal"
$ws.Range("G47").Value = "'1"
$ws.Range("G47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'0.002890"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'0.83%"
$ws.Range("E48").Style = "Normal"
$ws.Range("G48").Value = "'1"
$ws.Range("G48").Style = "Normal"

# Row 49
$ws.Range("D49").Value = "'0.001690"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-6.29%"
$ws.Range("E49").Style = "Normal"
$ws.Range("G49").Value = "'1"
$ws.Range("G49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'0.00002105"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.88%"
$ws.Range("E50").Style = "Normal"
$ws.Range("G50").Value = "'1"
$ws.Range("G50").Style = "Normal"

# Row 51
$ws.Range("D51").Value = "'0.0002005"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.88%"
$ws.Range("E51").Style = "Normal"
$ws.Range("G51").Value = "'1"
$ws.Range("G51").Style = "Normal"

